$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Commission" header
$ws.Range("G1").Value = "Commission"

# New data rows (Date, Type, Stock, Price, Lot, Total, Commission)
$data = @(
    @("26-10-2022", "Buy", "EREGL.IS", 30.28, 10, 302.8, 0.63),
    @("27-10-2022", "Buy", "SISE.IS", 33.92, 10, 339.2, 0.71),
    @("27-10-2022", "Buy", "EREGL.IS", 30.86, 15, 462.9, 0.97),
    @("07-11-2022", "Buy", "SISE.IS", 34.2, 20, 684, 1.43),
    @("08-11-2022", "Buy", "EREGL.IS", 36, 10, 360, 0.75),
    @("08-11-2022", "Buy", "THYAO.IS", 107.7, 3, 323.1, 0.68),
    @("10-11-2022", "Buy", "TUKAS.IS", 20.32, 25, 508, 1.06)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $row++
}
